$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Goal (per commit "feat: add 2022-Q1 data"):
#   * Insert a new sheet "2022-Q1" between "2021-Q4" and "总计", with the same
#     per-fund holdings layout as the other quarter sheets.
#   * Add a matching summary row to the "总计" (grand total) sheet, on top of
#     (above) the existing rows, shifting the older quarters down.
# ---------------------------------------------------------------------------

$q3 = $wb.Worksheets.Item(1)
$q4 = $wb.Worksheets.Item(2)
$totalOld = $wb.Worksheets.Item(3)

# Remove the original "总计" sheet first so the internal sheet-id counter is
# free to hand sheetId=3 to the new quarter sheet (matches target layout:
# 2022-Q1 -> id 3, 总计 -> id 4, recreated afterwards).
$totalOld.Delete()

# ---------------------------------------------------------------------------
# 1) New "2022-Q1" sheet - duplicate "2021-Q4" so header style/borders and
#    page setup (margins etc.) come along for free, then overwrite the data.
# ---------------------------------------------------------------------------
$q4.Copy($null, $q4)
$q1 = $wb.Worksheets.Item(3)
$q1.Name = "2022-Q1"

$q1.Range("D2").Value = "'0.29"
$q1.Range("D2").Style = "Normal"
$q1.Range("E2").Value = "'94.14"
$q1.Range("E2").Style = "Normal"
$q1.Range("F2").Value = "'2.54"
$q1.Range("F2").Style = "Normal"
$q1.Range("G2").Value = "'0.0074"
$q1.Range("G2").Style = "Normal"
$q1.Range("H2").Value = 4

# ---------------------------------------------------------------------------
# 2) Rebuild "总计" sheet at the end. Duplicating "2021-Q3" purely to inherit
#    its sheetPr/page-setup skeleton, then wipe the borrowed grid and
#    repopulate header + the 3 summary rows (new quarter on top).
# ---------------------------------------------------------------------------
$q3.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$total = $wb.Worksheets.Item($wb.Worksheets.Count)
$total.Name = "总计"

$total.Range("E1:H2").EntireColumn.Delete()
$total.Rows.Item(1).Delete()
$total.Rows.Item(1).Delete()

$total.PageSetup.LeftMargin = 54
$total.PageSetup.RightMargin = 54
$total.PageSetup.TopMargin = 72
$total.PageSetup.BottomMargin = 72
$total.PageSetup.HeaderMargin = 36
$total.PageSetup.FooterMargin = 36

$q4.Range("B1:D1").Copy($total.Range("B1:D1"))
$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$q4.Range("A2").Copy($total.Range("A2"))
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.01

$q4.Range("A2").Copy($total.Range("A3"))
$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 1
$total.Range("D3").Value = 0.01

$q4.Range("A2").Copy($total.Range("A4"))
$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q3"
$total.Range("C4").Value = 1
$total.Range("D4").Value = 0.01

Write-Host "done"
